# The edit rotates the species-observation data among rows 6-10:
#   new row 6  <- old row 7
#   new row 7  <- old row 10
#   new row 8  <- old row 6
#   new row 9  <- old row 8   (and gains the stray empty "L" cell old row 8 had)
#   new row 10 <- old row 9
#
# Columns A,B,D,E,F,G,H,Q,R are the ones that actually carry different data
# from row to row; every other column (C, I, J, K, N, P, S, T, U, V, W, Y, Z,
# AA, AB, AD, AE, AF, AG, AT, AW, AX, AY, ...) is identical across rows 6-10
# so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    6  = @{ A = 111671159; B = 81248;  D = "NT"; E = 1312;   F = "Gammelgransskål"; G = "Pseudographis pinicola";          H = "(Nyl.) Rehm" ;                    Q = 558006.0394731871; R = 7067389.087574247 }
    7  = @{ A = 111671165; B = 78578;  D = "NT"; E = 6458;   F = "Lunglav";          G = "Lobaria pulmonaria";              H = "(L.) Hoffm.";                     Q = 558014.2710882163; R = 7067448.175823289 }
    8  = @{ A = 111671148; B = 89405;  D = "NT"; E = 1202;   F = "Ullticka";         G = "Phellinidium ferrugineofuscum";   H = "(P.Karst.) Fiasson & Niemelä";    Q = 557930.9937661786; R = 7067802.902090888 }
    9  = @{ A = 111670690; B = 96348;  D = "VU"; E = 220787; F = "Knärot";           G = "Goodyera repens";                 H = "(L.) R. Br.";                     Q = 557809.1117697239; R = 7067699.199123298 }
    10 = @{ A = 111670912; B = 78578;  D = "NT"; E = 6458;   F = "Lunglav";          G = "Lobaria pulmonaria";              H = "(L.) Hoffm.";                     Q = 557803.3534448177; R = 7067771.317107533 }
}

# The stray, content-less "L" cell that currently sits on row 8 travels
# along with that row's data, landing on row 9 (old row 8 -> new row 9).
# Copy the blank cell first (before the row 8 data itself is overwritten
# below) and then remove it from its old spot.
$ws.Range("L8").Copy($ws.Range("L9"))
$ws.Range("L8").ClearContents()

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value2  = $vals.A   # A - Id
    $ws.Cells.Item($r, 2).Value2  = $vals.B   # B - Taxonsorteringsordning
    $ws.Cells.Item($r, 4).Value2  = $vals.D   # D - Rödlistade
    $ws.Cells.Item($r, 5).Value2  = $vals.E   # E - TaxonId
    $ws.Cells.Item($r, 6).Value2  = $vals.F   # F - Artnamn
    $ws.Cells.Item($r, 7).Value2  = $vals.G   # G - Vetenskapligt namn
    $ws.Cells.Item($r, 8).Value2  = $vals.H   # H - Auktor
    $ws.Cells.Item($r, 17).Value2 = $vals.Q   # Q - Ost
    $ws.Cells.Item($r, 18).Value2 = $vals.R   # R - Nord
}
